$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy formatting from row 11 to the new row 12 so the new row matches
# the look of the existing "26-09-2025" rows (style ids 3 and 4).
$ws.Range("A11:B11").Copy() | Out-Null
$ws.Range("A12:B12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Set the literal text values (these already exist in the shared string
# table, so they will reuse the same shared-string entries as row 11).
$ws.Range("A12").Value = "26-09-2025"
$ws.Range("B12").Value = "The price of gold in India today is ₹11,488 per gram for 24 karat gold, ₹10,530 per gram for 22 karat gold and ₹8,616 per gram for 18 karat gold (also called 999 gold)."
